# Generate Report for Archive
#
# 1) Update the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F2/E3:F3 and the
#    "Status" column (C2:C3) on each per-locale sheet).
# 2) Narrow the "status" columns' widths (Overview columns E & F, and
#    column C on the zh-cn / de-de sheets) to match the new narrower
#    "In Translation" header footprint used in the archived report.

$wb = $excel.ActiveWorkbook

# --- 1. Text update -------------------------------------------------------

$overview = $wb.Worksheets.Item(1)
$overview.Range("E2:F2").Value = "In Translation"
$overview.Range("E3:F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item(2)
$zhcn.Range("C2:C3").Value = "In Translation"

$dede = $wb.Worksheets.Item(3)
$dede.Range("C2:C3").Value = "In Translation"

# --- 2. Column width update -------------------------------------------------
# ColumnWidth is expressed in characters and gets quantized to whole pixels
# (pixels = floor(width*6 + 0.5) + 5, stored width = pixels/6), same as real
# Excel. 12.5 is the input that lands on the pixel bucket closest to the
# target stored width (~13.41 chars) for these columns.

$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
